$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 42 (shifts existing rows 42-100 down to 43-101,
# and Excel's dimension / formatting follow automatically, matching the target diff).
$ws.Rows(42).Insert()

$ws.Cells.Item(42, 1).Value  = 8
$ws.Cells.Item(42, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42, 3).Value  = "Coquimbo"
$ws.Cells.Item(42, 4).Value  = 44915
$ws.Cells.Item(42, 5).Value  = 4
$ws.Cells.Item(42, 6).Value  = 100112030
$ws.Cells.Item(42, 7).Value  = "Poroto granado"
$ws.Cells.Item(42, 8).Value  = "Sin especificar"
$ws.Cells.Item(42, 9).Value  = "Primera"
$ws.Cells.Item(42, 10).Value = 520
$ws.Cells.Item(42, 11).Value = 35000
$ws.Cells.Item(42, 12).Value = 36000
$ws.Cells.Item(42, 13).Value = 35500
$ws.Cells.Item(42, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 16).Value = 1420
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
